# Update session 42 documents.
#
# 1) Shade the last 4 data rows (Mon Dec 4 .. Mon Dec 11 / FINAL EXAM) of the
#    main calendar table with the same "gray 50%, background 1, darker 50%"
#    cell shading already used on the earlier weeks of the table.
# 2) The _GoBack bookmark (Word's "last edit location" marker) moves from the
#    blank paragraph above "Archive:" to the end of the "Archive:" run, which
#    is simply a side effect of where the edits above were made.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# ---------------------------------------------------------------------
# Step 1: add shading to the 16 table cells (4 rows x 4 columns) that
# cover the "Mon, Dec 4" ... "FINAL EXAM" rows of the first table.
# ---------------------------------------------------------------------
$rowStartNeedle = "Mon, Dec 4"
$rowStartIdx = $xml.IndexOf($rowStartNeedle)
$trOpenIdx = $xml.LastIndexOf("<w:tr ", $rowStartIdx)

$rowEndNeedle = "FINAL EXAM"
$rowEndIdx = $xml.IndexOf($rowEndNeedle, $rowStartIdx)
$trCloseTag = "</w:tr>"
$trCloseIdx = $xml.IndexOf($trCloseTag, $rowEndIdx) + $trCloseTag.Length

$segment = $xml.Substring($trOpenIdx, $trCloseIdx - $trOpenIdx)

$shdXml = '<w:shd w:val="clear" w:color="auto" w:fill="808080" w:themeFill="background1" w:themeFillShade="80"/>'

$colWidths = @("1254", "4321", "1707", "1708")
foreach ($w in $colWidths) {
    $oldTcPr = '<w:tcW w:w="' + $w + '" w:type="dxa"/></w:tcPr>'
    $newTcPr = '<w:tcW w:w="' + $w + '" w:type="dxa"/>' + $shdXml + '</w:tcPr>'
    $segment = $segment.Replace($oldTcPr, $newTcPr)
}

$xml = $xml.Substring(0, $trOpenIdx) + $segment + $xml.Substring($trCloseIdx)

# ---------------------------------------------------------------------
# Step 2: move the _GoBack bookmark from the empty paragraph to right
# after the "Archive:" run.
# ---------------------------------------------------------------------
$bookmarkXml = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$xml = $xml.Replace($bookmarkXml, "")

$archiveOld = '<w:t>Archive:</w:t></w:r></w:p>'
$archiveNew = '<w:t>Archive:</w:t></w:r>' + $bookmarkXml + '</w:p>'
$xml = $xml.Replace($archiveOld, $archiveNew)

$d.WordOpenXML = $xml
